$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 41; this pushes the existing rows 41-73 down
# to 42-74 (same as the weekly roll in the source diff), keeping every
# column value intact for the shifted rows.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with this week's record.
$ws.Range("A41").Value2 = 2
$ws.Range("B41").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C41").Value2 = "Coquimbo"
$ws.Range("D41").Value2 = 44818
$ws.Range("E41").Value2 = 4
$ws.Range("F41").Value2 = 100112022
$ws.Range("G41").Value2 = "Arveja Verde"
$ws.Range("H41").Value2 = "Perfection"
$ws.Range("I41").Value2 = "Primera"
$ws.Range("J41").Value2 = 400
$ws.Range("K41").Value2 = 24000
$ws.Range("L41").Value2 = 27000
$ws.Range("M41").Value2 = 25500
$ws.Range("N41").Value2 = "`$/malla 25 kilos"
$ws.Range("O41").Value2 = "Provincia de Limarí"
$ws.Range("P41").Value2 = 1020
$ws.Range("Q41").Value2 = 25
$ws.Range("R41").Value2 = "Hortaliza"

# Make sure the date column keeps the date number format used by the
# rest of the column (same style as D42:D74).
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
